$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format D/E columns as Text so numeric-looking strings (e.g. "1.00", "486.42")
# are preserved exactly as text instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '55.396.27'
$ws.Range('E2').Value = '  -4.89%  '
$ws.Range('D3').Value = '2.891.62'
$ws.Range('E3').Value = '  -5.56%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '486.42'
$ws.Range('E5').Value = '  -6.84%  '
$ws.Range('D6').Value = '131.56'
$ws.Range('E6').Value = '  -7.66%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '0.416'
$ws.Range('E8').Value = '  -7.06%  '
$ws.Range('D9').Value = '7.10'
$ws.Range('E9').Value = '  -5.42%  '
$ws.Range('D10').Value = '0.103'
$ws.Range('E10').Value = '  -8.33%  '
$ws.Range('D11').Value = '0.344'
$ws.Range('E11').Value = '  -7.00%  '
$ws.Range('D12').Value = '3.372.33'
$ws.Range('E12').Value = '  -8.68%  '
$ws.Range('E13').Value = '  -4.63%  '
$ws.Range('D14').Value = '25.51'
$ws.Range('E14').Value = '  -4.59%  '
$ws.Range('D15').Value = '0.0000156'
$ws.Range('E15').Value = '  -9.01%  '
$ws.Range('D16').Value = '55.370.80'
$ws.Range('E16').Value = '  -4.88%  '
$ws.Range('D17').Value = '5.93'
$ws.Range('E17').Value = '  -4.75%  '
$ws.Range('D18').Value = '2.897.19'
$ws.Range('E18').Value = '  -5.67%  '
$ws.Range('D19').Value = '12.31'
$ws.Range('E19').Value = '  -5.97%  '
$ws.Range('D20').Value = '7.57'
$ws.Range('E20').Value = '  -7.15%  '
$ws.Range('D21').Value = '310.84'
$ws.Range('E21').Value = '  -8.21%  '
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D22').Value = '5.79'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '0.474'
$ws.Range('E24').Value = '  -5.73%  '
$ws.Range('D25').Value = '61.58'
$ws.Range('E25').Value = '  -5.89%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = '0.159'
$ws.Range('E27').Value = '  -6.50%  '
$ws.Range('D28').Value = '0.0₃0830'
$ws.Range('E28').Value = '  -13.79%  '
$ws.Range('D29').Value = '6.31'
$ws.Range('E29').Value = '  -9.22%  '
$ws.Range('D30').Value = '6.91'
$ws.Range('E30').Value = '  -8.06%  '
$ws.Range('D31').Value = '1.72'
$ws.Range('E31').Value = '  -6.67%  '
$ws.Range('D32').Value = '19.51'
$ws.Range('E32').Value = '  -7.85%  '
$ws.Range('D33').Value = '1.11'
$ws.Range('E33').Value = '  -9.54%  '
$ws.Range('D34').Value = '148.72'
$ws.Range('E34').Value = '  -5.65%  '
$ws.Range('D35').Value = '4.34'
$ws.Range('E35').Value = '  -9.59%  '
$ws.Range('D36').Value = '5.54'
$ws.Range('E36').Value = '  -7.08%  '
$ws.Range('D37').Value = '24.05'
$ws.Range('E37').Value = '  -5.64%  '
$ws.Range('D38').Value = '1.18'
$ws.Range('E38').Value = '  -9.88%  '
$ws.Range('D39').Value = '0.0647'
$ws.Range('E39').Value = '  -6.87%  '
$ws.Range('D40').Value = '2.921.88'
$ws.Range('E40').Value = '  -5.77%  '
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').Value = '35.79'
$ws.Range('E42').Value = '  -5.37%  '
$ws.Range('E43').Value = '  -7.95%  '
$ws.Range('D44').Value = '0.624'
$ws.Range('E44').Value = '  -6.58%  '
$ws.Range('D45').Value = '2.082.47'
$ws.Range('E45').Value = '  -10.87%  '
$ws.Range('D46').Value = '1.31'
$ws.Range('E46').Value = '  -10.10%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').Value = '5.77'
$ws.Range('E47').Value = '  -5.02%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '0.909'
$ws.Range('E48').Value = '  -10.67%  '
$ws.Range('D49').Value = '0.0225'
$ws.Range('E49').Value = '  -6.92%  '
$ws.Range('D50').Value = '18.51'
$ws.Range('E50').Value = '  -6.85%  '
$ws.Range('D51').Value = '0.0829'
$ws.Range('E51').Value = '  -8.07%  '

# Restore the original (default) cell style on the Text-formatted range so the
# workbook styling matches the source formatting (no lingering "@" format).
$ws.Range("D2:E51").Style = "Normal"
